# Auto-generated edit script applying scheduled profit-recalculation updates
# to the Zodiark_Profits workbook (columns H-N: pricing/profit figures) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 533.4167
$ws.Range("I41").Value = 285.66666
$ws.Range("J41").Value = 1276.6666
$ws.Range("K41").Value = 285.66666
$ws.Range("L41").Value = 1276.6666
$ws.Range("M41").Value = 154.33334
$ws.Range("N41").Value = -2156.6666
# Row 98
$ws.Range("H98").Value = 6166.769
$ws.Range("I98").Value = 2709.75
$ws.Range("J98").Value = 11698
$ws.Range("K98").Value = 2709.75
$ws.Range("L98").Value = 11698
$ws.Range("M98").Value = -1211.75
$ws.Range("N98").Value = -14694
# Row 122
$ws.Range("H122").Value = 6166.769
$ws.Range("I122").Value = 2709.75
$ws.Range("J122").Value = 11698
$ws.Range("K122").Value = 8129.25
$ws.Range("L122").Value = 35094
$ws.Range("M122").Value = -5679.25
$ws.Range("N122").Value = -39994
# Row 125
$ws.Range("H125").Value = 1853.5714
$ws.Range("J125").Value = 1799.6666
$ws.Range("L125").Value = 16196.9994
$ws.Range("N125").Value = -21116.9994
# Row 138
$ws.Range("H138").Value = 2074.9092
$ws.Range("J138").Value = 2380.2273
$ws.Range("L138").Value = 7140.6819
$ws.Range("N138").Value = -17420.6819

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1461.9891
$ws.Range("I32").Value = 1096.966
$ws.Range("K32").Value = 1096.966
$ws.Range("M32").Value = -809.9659999999999
# Row 45
$ws.Range("H45").Value = 2323.8333
$ws.Range("I45").Value = 2482.25
$ws.Range("J45").Value = 2007
$ws.Range("K45").Value = 2482.25
$ws.Range("L45").Value = 2007
$ws.Range("M45").Value = -2105.25
$ws.Range("N45").Value = -2761
# Row 49
$ws.Range("H49").Value = 69949
$ws.Range("J49").Value = 69949
$ws.Range("L49").Value = 69949
$ws.Range("N49").Value = -70469
# Row 74
$ws.Range("H74").Value = 1463.2325
$ws.Range("I74").Value = 1490.826
$ws.Range("J74").Value = 1431.5
$ws.Range("K74").Value = 1490.826
$ws.Range("L74").Value = 1431.5
$ws.Range("M74").Value = -616.826
$ws.Range("N74").Value = -3179.5
# Row 77
$ws.Range("H77").Value = 1463.2325
$ws.Range("I77").Value = 1490.826
$ws.Range("J77").Value = 1431.5
$ws.Range("K77").Value = 7454.13
$ws.Range("L77").Value = 7157.5
$ws.Range("M77").Value = -3086.13
$ws.Range("N77").Value = -15893.5
# Row 122
$ws.Range("H122").Value = 2848.6072
$ws.Range("I122").Value = 2450.52
$ws.Range("K122").Value = 7351.559999999999
$ws.Range("M122").Value = -4901.559999999999
# Row 132
$ws.Range("H132").Value = 6011.615
$ws.Range("I132").Value = 7955.857
$ws.Range("K132").Value = 23867.571
$ws.Range("M132").Value = -21337.571

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 900.1818
$ws.Range("J20").Value = 899.2857
$ws.Range("L20").Value = 899.2857
$ws.Range("N20").Value = -1393.2857
# Row 94
$ws.Range("H94").Value = 953.6539
$ws.Range("I94").Value = 991.8261
$ws.Range("K94").Value = 991.8261
$ws.Range("M94").Value = -540.8261
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 134
$ws.Range("H134").Value = 17217.666
$ws.Range("I134").Value = 23022.75
$ws.Range("J134").Value = 11412.583
$ws.Range("K134").Value = 69068.25
$ws.Range("L134").Value = 34237.749
$ws.Range("M134").Value = -66533.25
$ws.Range("N134").Value = -39307.749

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2651.5
$ws.Range("I22").Value = 2592.7693
$ws.Range("K22").Value = 2592.7693
$ws.Range("M22").Value = -2242.7693
# Row 31
$ws.Range("H31").Value = 4360.8716
$ws.Range("J31").Value = 6030.3076
$ws.Range("L31").Value = 6030.3076
$ws.Range("N31").Value = -6620.3076
# Row 34
$ws.Range("H34").Value = 4360.8716
$ws.Range("J34").Value = 6030.3076
$ws.Range("L34").Value = 6030.3076
$ws.Range("N34").Value = -6434.3076
# Row 74
$ws.Range("H74").Value = 34576.832
$ws.Range("J74").Value = 36678.2
$ws.Range("L74").Value = 36678.2
$ws.Range("N74").Value = -38426.2
# Row 77
$ws.Range("H77").Value = 34576.832
$ws.Range("J77").Value = 36678.2
$ws.Range("L77").Value = 110034.6
$ws.Range("N77").Value = -118770.6

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 3777.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 98
$ws.Range("H98").Value = 509.5
$ws.Range("J98").Value = 567
$ws.Range("L98").Value = 1701
$ws.Range("N98").Value = -4697
# Row 107
$ws.Range("H107").Value = 932.6667
$ws.Range("J107").Value = 1050.2142
$ws.Range("L107").Value = 3150.6426
$ws.Range("N107").Value = -6990.642599999999
# Row 133
$ws.Range("H133").Value = 8000
$ws.Range("I133").Value = 1000
$ws.Range("J133").Value = 11500
$ws.Range("K133").Value = 3000
$ws.Range("L133").Value = 34500
$ws.Range("M133").Value = 2060
$ws.Range("N133").Value = -44620

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 42980
$ws.Range("I52").Value = 40000
$ws.Range("J52").Value = 44966.668
$ws.Range("K52").Value = 40000
$ws.Range("L52").Value = 44966.668
$ws.Range("M52").Value = -39741
$ws.Range("N52").Value = -45484.668
# Row 70
$ws.Range("H70").Value = 55567.332
$ws.Range("I70").Value = 64756.555
$ws.Range("K70").Value = 64756.555
$ws.Range("M70").Value = -64486.555
# Row 73
$ws.Range("H73").Value = 55567.332
$ws.Range("I73").Value = 64756.555
$ws.Range("K73").Value = 64756.555
$ws.Range("M73").Value = -63820.555
# Row 123
$ws.Range("H123").Value = 102908
$ws.Range("J123").Value = 102908
$ws.Range("L123").Value = 102908
$ws.Range("N123").Value = -107808
# Row 126
$ws.Range("H126").Value = 10269630
$ws.Range("I126").Value = 6687.2
$ws.Range("K126").Value = 20061.6
$ws.Range("M126").Value = -17591.6
# Row 132
$ws.Range("H132").Value = 9746.046
$ws.Range("I132").Value = 9412.647000000001
$ws.Range("K132").Value = 28237.941
$ws.Range("M132").Value = -25707.941

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4453.6523
$ws.Range("I40").Value = 4621.45
$ws.Range("K40").Value = 4621.45
$ws.Range("M40").Value = -4485.45
# Row 61
$ws.Range("H61").Value = 5848.7
$ws.Range("I61").Value = 5864.7617
$ws.Range("J61").Value = 5811.222
$ws.Range("K61").Value = 5864.7617
$ws.Range("L61").Value = 5811.222
$ws.Range("M61").Value = -5662.7617
$ws.Range("N61").Value = -6215.222
# Row 68
$ws.Range("H68").Value = 5762.091
$ws.Range("J68").Value = 11594.667
$ws.Range("L68").Value = 11594.667
$ws.Range("N68").Value = -13092.667
# Row 71
$ws.Range("H71").Value = 5762.091
$ws.Range("J71").Value = 11594.667
$ws.Range("L71").Value = 57973.335
$ws.Range("N71").Value = -65461.335
# Row 113
$ws.Range("H113").Value = 5848.7
$ws.Range("I113").Value = 5864.7617
$ws.Range("J113").Value = 5811.222
$ws.Range("K113").Value = 5864.7617
$ws.Range("L113").Value = 5811.222
$ws.Range("M113").Value = -3694.7617
$ws.Range("N113").Value = -10151.222
# Row 122
$ws.Range("H122").Value = 4898.2
$ws.Range("I122").Value = 4932.875
$ws.Range("K122").Value = 14798.625
$ws.Range("M122").Value = -12348.625
# Row 132
$ws.Range("H132").Value = 10325
$ws.Range("I132").Value = 9766.666999999999
$ws.Range("K132").Value = 29300.001
$ws.Range("M132").Value = -26770.001

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 6396.2
$ws.Range("J2").Value = 6396.2
$ws.Range("L2").Value = 6396.2
$ws.Range("N2").Value = -6620.2
# Row 54
$ws.Range("H54").Value = 28993.4
$ws.Range("J54").Value = 28993.4
$ws.Range("L54").Value = 28993.4
$ws.Range("N54").Value = -30033.4
# Row 81
$ws.Range("H81").Value = 4739.75
$ws.Range("J81").Value = 7479.5
$ws.Range("L81").Value = 14959
$ws.Range("N81").Value = -17081
# Row 84
$ws.Range("H84").Value = 4739.75
$ws.Range("J84").Value = 7479.5
$ws.Range("L84").Value = 74795
$ws.Range("N84").Value = -85403
# Row 122
$ws.Range("H122").Value = 15154924
$ws.Range("I122").Value = 16132306
$ws.Range("K122").Value = 48396918
$ws.Range("M122").Value = -48394468
# Row 132
$ws.Range("H132").Value = 4231.5
$ws.Range("I132").Value = 4231.5
$ws.Range("K132").Value = 12694.5
$ws.Range("M132").Value = -10164.5

